# Insert 3 new weekly-report rows above the current row 586, pushing the
# existing data (rows 586-686) down to rows 589-689, then fill the three
# newly-opened rows (586-588) with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 586..686 down by 3 rows (inserts 3 blank rows at 586..588).
$ws.Range("A586:T588").EntireRow.Insert()

# Constant values shared by every data row of this sub-sheet.
$mercadoId   = 2
$mercado     = "Comercializadora del Agro de Limarí"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100101
$producto    = "Berries"
$categoriaId = 100112025
$categoria   = "Frutilla"
$variedad    = "Sin especificar"
$unidad      = "$/bandeja 7 kilos"
$origen      = "Provincia de Melipilla"
$kgUnidad    = 7

# New week's data (fecha serial 45218 = 2023-10-19) for the three quality grades.
$nuevasFilas = @(
    @{ Fila = 586; Calidad = "Especial"; Volumen = 360; PMin = 14000; PMax = 15000; PProm = 14500; PKg = 2071 },
    @{ Fila = 587; Calidad = "Primera";  Volumen = 400; PMin = 11000; PMax = 12000; PProm = 11500; PKg = 1643 },
    @{ Fila = 588; Calidad = "Segunda";  Volumen = 320; PMin = 8000;  PMax = 9000;  PProm = 8500;  PKg = 1214 }
)

foreach ($fila in $nuevasFilas) {
    $r = $fila.Fila
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = 45218
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $fila.Calidad
    $ws.Cells.Item($r, 13).Value = $fila.Volumen
    $ws.Cells.Item($r, 14).Value = $fila.PMin
    $ws.Cells.Item($r, 15).Value = $fila.PMax
    $ws.Cells.Item($r, 16).Value = $fila.PProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $fila.PKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
